$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 72338.64
$ws.Range("I2").Value = 460.33334
$ws.Range("J2").Value = 126247.375
$ws.Range("K2").Value = 460.33334
$ws.Range("L2").Value = 126247.375
$ws.Range("M2").Value = -347.33334
$ws.Range("N2").Value = -126473.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1224.5714
$ws.Range("I12").Value = 417.5
$ws.Range("K12").Value = 417.5
$ws.Range("M12").Value = -247.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3676.318
$ws.Range("J19").Value = 5404.727
$ws.Range("L19").Value = 5404.727
$ws.Range("N19").Value = -5754.727

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3945.182
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 8224.25
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 8224.25
$ws.Range("M40").Value = -1325
$ws.Range("N40").Value = -8574.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1735.3077
$ws.Range("I127").Value = 1978.2727
$ws.Range("K127").Value = 5934.8181
$ws.Range("M127").Value = -974.8181000000004

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1891
$ws.Range("J138").Value = 2392.3333
$ws.Range("L138").Value = 7176.999899999999
$ws.Range("N138").Value = -17456.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 5221.6875
$ws.Range("J141").Value = 4748.25
$ws.Range("L141").Value = 14244.75
$ws.Range("N141").Value = -24604.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 2466
$ws.Range("I4").Value = 1811.7916
$ws.Range("J4").Value = 5082.8335
$ws.Range("K4").Value = 1811.7916
$ws.Range("L4").Value = 5082.8335
$ws.Range("M4").Value = -1695.7916
$ws.Range("N4").Value = -5314.8335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7085.65
$ws.Range("I45").Value = 8240.214
$ws.Range("K45").Value = 8240.214
$ws.Range("M45").Value = -7863.214

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3715.6738
$ws.Range("I61").Value = 2775.75
$ws.Range("K61").Value = 2775.75
$ws.Range("M61").Value = -2563.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1554.4736
$ws.Range("I122").Value = 1410.9286
$ws.Range("K122").Value = 4232.7858
$ws.Range("M122").Value = -1782.7858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3715.6738
$ws.Range("I136").Value = 2775.75
$ws.Range("K136").Value = 8327.25
$ws.Range("M136").Value = -5777.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1687.1875
$ws.Range("I31").Value = 1035.909
$ws.Range("J31").Value = 3120
$ws.Range("K31").Value = 1035.909
$ws.Range("L31").Value = 3120
$ws.Range("M31").Value = -740.9090000000001
$ws.Range("N31").Value = -3710

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1687.1875
$ws.Range("I34").Value = 1035.909
$ws.Range("J34").Value = 3120
$ws.Range("K34").Value = 1035.909
$ws.Range("L34").Value = 3120
$ws.Range("M34").Value = -833.9090000000001
$ws.Range("N34").Value = -3524

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1267.3704
$ws.Range("I107").Value = 1008.8421
$ws.Range("J107").Value = 1881.375
$ws.Range("K107").Value = 1008.8421
$ws.Range("L107").Value = 1881.375
$ws.Range("M107").Value = 911.1579
$ws.Range("N107").Value = -5721.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 11352.7
$ws.Range("I132").Value = 11259
$ws.Range("J132").Value = 11571.333
$ws.Range("K132").Value = 33777
$ws.Range("L132").Value = 34713.999
$ws.Range("M132").Value = -31247
$ws.Range("N132").Value = -39773.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 79999.5
$ws.Range("J138").Value = 79999.5
$ws.Range("L138").Value = 79999.5
$ws.Range("N138").Value = -90279.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9117.362999999999
$ws.Range("I3").Value = 7695.3335
$ws.Range("J3").Value = 15516.5
$ws.Range("K3").Value = 23086.0005
$ws.Range("L3").Value = 46549.5
$ws.Range("M3").Value = -22974.0005
$ws.Range("N3").Value = -46773.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3589.75
$ws.Range("I34").Value = 429.75
$ws.Range("K34").Value = 1289.25
$ws.Range("M34").Value = -1205.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 71607610
$ws.Range("I55").Value = 586647.75
$ws.Range("J55").Value = 100016000
$ws.Range("K55").Value = 1759943.25
$ws.Range("L55").Value = 300048000
$ws.Range("M55").Value = -1759766.25
$ws.Range("N55").Value = -300048354

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1752.909
$ws.Range("I86").Value = 464.33334
$ws.Range("J86").Value = 3299.2
$ws.Range("K86").Value = 1393.00002
$ws.Range("L86").Value = 9897.599999999999
$ws.Range("M86").Value = -207.0000199999999
$ws.Range("N86").Value = -12269.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 8147.4
$ws.Range("I87").Value = 8219.333000000001
$ws.Range("J87").Value = 7500
$ws.Range("K87").Value = 24657.999
$ws.Range("L87").Value = 22500
$ws.Range("M87").Value = -23409.999
$ws.Range("N87").Value = -24996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 1752.909
$ws.Range("I89").Value = 464.33334
$ws.Range("J89").Value = 3299.2
$ws.Range("K89").Value = 4179.00006
$ws.Range("L89").Value = 29692.8
$ws.Range("M89").Value = 1748.99994
$ws.Range("N89").Value = -41548.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 8147.4
$ws.Range("I90").Value = 8219.333000000001
$ws.Range("J90").Value = 7500
$ws.Range("K90").Value = 73973.997
$ws.Range("L90").Value = 67500
$ws.Range("M90").Value = -67733.997
$ws.Range("N90").Value = -79980

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2643.4614
$ws.Range("I107").Value = 702.5
$ws.Range("J107").Value = 3506.111
$ws.Range("K107").Value = 2107.5
$ws.Range("L107").Value = 10518.333
$ws.Range("M107").Value = -187.5
$ws.Range("N107").Value = -14358.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 6690
$ws.Range("I108").Value = 510
$ws.Range("K108").Value = 1530
$ws.Range("M108").Value = 1350

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 7678.12
$ws.Range("I109").Value = 2530.3
$ws.Range("J109").Value = 11110
$ws.Range("K109").Value = 7590.900000000001
$ws.Range("L109").Value = 33330
$ws.Range("M109").Value = -6550.900000000001
$ws.Range("N109").Value = -35410

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 62502916
$ws.Range("I113").Value = 410
$ws.Range("J113").Value = 76926570
$ws.Range("K113").Value = 1230
$ws.Range("L113").Value = 230779710
$ws.Range("M113").Value = 940
$ws.Range("N113").Value = -230784050

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1395
$ws.Range("I114").Value = 989
$ws.Range("J114").Value = 1496.5
$ws.Range("K114").Value = 2967
$ws.Range("L114").Value = 4489.5
$ws.Range("M114").Value = 287
$ws.Range("N114").Value = -10997.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 6510.3335
$ws.Range("I129").Value = 1074.125
$ws.Range("J129").Value = 50000
$ws.Range("K129").Value = 3222.375
$ws.Range("L129").Value = 150000
$ws.Range("M129").Value = 1777.625
$ws.Range("N129").Value = -160000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4125.75
$ws.Range("I132").Value = 4128
$ws.Range("K132").Value = 12384
$ws.Range("M132").Value = -9854

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 61001.715
$ws.Range("I7").Value = 70335.336
$ws.Range("K7").Value = 70335.336
$ws.Range("M7").Value = -70223.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 5000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 5000
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -4830

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3323.2942
$ws.Range("I22").Value = 1539.091
$ws.Range("J22").Value = 6594.3335
$ws.Range("K22").Value = 1539.091
$ws.Range("L22").Value = 6594.3335
$ws.Range("M22").Value = -1244.091
$ws.Range("N22").Value = -7184.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 3323.2942
$ws.Range("I27").Value = 1539.091
$ws.Range("J27").Value = 6594.3335
$ws.Range("K27").Value = 1539.091
$ws.Range("L27").Value = 6594.3335
$ws.Range("M27").Value = -1432.091
$ws.Range("N27").Value = -6808.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4348
$ws.Range("J46").Value = 4348
$ws.Range("L46").Value = 4348
$ws.Range("N46").Value = -4724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1546.0344
$ws.Range("I55").Value = 323.29413
$ws.Range("J55").Value = 3278.25
$ws.Range("K55").Value = 323.29413
$ws.Range("L55").Value = 3278.25
$ws.Range("M55").Value = -150.29413
$ws.Range("N55").Value = -3624.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 19004.084
$ws.Range("I61").Value = 16714.4
$ws.Range("K61").Value = 16714.4
$ws.Range("M61").Value = -16512.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2565.92
$ws.Range("I100").Value = 2377.647
$ws.Range("J100").Value = 2966
$ws.Range("K100").Value = 2377.647
$ws.Range("L100").Value = 2966
$ws.Range("M100").Value = -1836.647
$ws.Range("N100").Value = -4048

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 19004.084
$ws.Range("I113").Value = 16714.4
$ws.Range("K113").Value = 16714.4
$ws.Range("M113").Value = -14544.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4396.72
$ws.Range("I122").Value = 3244.818
$ws.Range("K122").Value = 9734.454000000002
$ws.Range("M122").Value = -7284.454000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 61001.715
$ws.Range("I126").Value = 70335.336
$ws.Range("K126").Value = 211006.008
$ws.Range("M126").Value = -208536.008

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9590.916999999999
$ws.Range("J62").Value = 11399.111
$ws.Range("L62").Value = 11399.111
$ws.Range("N62").Value = -12647.111

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 9590.916999999999
$ws.Range("J65").Value = 11399.111
$ws.Range("L65").Value = 56995.55500000001
$ws.Range("N65").Value = -63235.55500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5167.5
$ws.Range("I136").Value = 5838.6665
$ws.Range("J136").Value = 3959.4
$ws.Range("K136").Value = 17515.9995
$ws.Range("L136").Value = 11878.2
$ws.Range("M136").Value = -14965.9995
$ws.Range("N136").Value = -16978.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 69043.17
$ws.Range("J139").Value = 69043.17
$ws.Range("L139").Value = 69043.17
$ws.Range("N139").Value = -79323.17

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 114749.836
$ws.Range("J141").Value = 114749.836
$ws.Range("L141").Value = 114749.836
$ws.Range("N141").Value = -125109.836
